$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the "Ator(es)" values for rows that now also include "Usuário Comum".
$ws.Range("C3").Value = "Administrador/Comprador/Organizadora/Usuário Comum"
$ws.Range("C4").Value = "Administrador/Comprador/Organizadora/Usuário Comum"
$ws.Range("C21").Value = "Usuário Comum"
$ws.Range("C31").Value = "Usuário Comum"

# Widen column C to fit the new, longer text (saved width rounds to 54).
$ws.Columns.Item(3).ColumnWidth = 53.17

# Move the active selection to E27.
$ws.Range("E27").Select()
